$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the 6 rows of old "Cloud Migration" leftover content (rows 708-713)
# so everything below shifts up by 6 and lines up with the new group boundaries.
$ws.Rows("708:713").Delete()

# Step 2: replace the Pub_id/Title of the (now 21-row) "Cloud Migration" group (rows 687-707)
# with the refreshed list of core publications about cloud migration.
$cloudMigration = @(
    @(687, "pub.1094505139", "CDOSim: Simulating Cloud Deployment Options for Software Migration Support"),
    @(688, "pub.1094808840", "An Extensible Architecture for Detecting Violations of a Cloud Environment&#x27;s Constraints During Legacy Software System Migration"),
    @(689, "pub.1119463839", "CloudGenius: Decision Support for Web Server Cloud Migration"),
    @(690, "pub.1093694541", "Cloudstep: A Step-by-Step Decision Process to Support Legacy Application Migration to the Cloud"),
    @(691, "pub.1094203398", "Software Engineering Challenges for Migration to the Service Cloud Paradigm"),
    @(692, "pub.1093422882", "Cloud Migration: A Case Study of Migrating an Enterprise IT System to IaaS"),
    @(693, "pub.1026275219", "How to adapt applications for the Cloud environment"),
    @(694, "pub.1094635108", "Migration of Multi-tier Applications to Infrastructure-as-a-Service Clouds: An Investigation Using Kernel-based Virtual Machines"),
    @(695, "pub.1095459235", "Legacy Application Migration to the Cloud: Practicability and Methodology"),
    @(696, "pub.1095784692", "Migration to Cloud as Real Option Investment decision under uncertainty"),
    @(697, "pub.1093590373", "Migrating Service-Oriented System to Cloud Computing: An Experience Report"),
    @(698, "pub.1015329134", "Cloud adoption"),
    @(699, "pub.1095092589", "Workload Migration into Clouds - Challenges, Experiences, Opportunities"),
    @(700, "pub.1031146575", "Application migration to cloud"),
    @(701, "pub.1029470968", "A tale of migration to cloud computing for sharing experiences and observations"),
    @(702, "pub.1094056660", "Size Estimation of Cloud Migration Projects with Cloud Migration Point (CMP)"),
    @(703, "pub.1095101569", "A Practical Architecture of Cloudification of Legacy Applications"),
    @(704, "pub.1017004215", "Automatic conformance checking for migrating software systems to cloud infrastructures and platforms"),
    @(705, "pub.1027646047", "The Cloud Adoption Toolkit: supporting cloud adoption decisions in the enterprise"),
    @(706, "pub.1063158434", "Cloudward bound"),
    @(707, "pub.1095280257", "Search-Based Genetic Optimization for Deployment and Reconfiguration of Software in the Cloud")
)

foreach ($rec in $cloudMigration) {
    $r = $rec[0]
    $ws.Cells.Item($r, 3).Value = $rec[1]
    $ws.Cells.Item($r, 4).Value = $rec[2]
}

Write-Host "done"